$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was inserted as row 200 ("Fruta / hortaliza, semanal"
# update), pushing the former rows 200-263 down to 201-264.
$ws.Rows.Item(200).EntireRow.Insert()

$ws.Range("A200").Value = 10
$ws.Range("B200").Value = "Vega Modelo de Temuco"
$ws.Range("C200").Value = "La Araucanía"
$ws.Range("D200").Value = 44588
$ws.Range("E200").Value = 9
$ws.Range("F200").Value = 100112009
$ws.Range("G200").Value = "Acelga"
$ws.Range("H200").Value = "Sin especificar"
$ws.Range("I200").Value = "Primera"
$ws.Range("J200").Value = 85
$ws.Range("K200").Value = 7000
$ws.Range("L200").Value = 8000
$ws.Range("M200").Value = 7588
$ws.Range("N200").Value = "`$/docena de atados (12 kilos)"
$ws.Range("O200").Value = "Provincia de Cautín"
$ws.Range("P200").Value = 632
$ws.Range("Q200").Value = 12
$ws.Range("R200").Value = "Hortaliza"
